$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.039.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.420.87'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '552.37'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.87'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.576'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.78'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.359'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '24.94'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.19%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.848.64'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.35%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.951.79'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.78%  '
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.420.36'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('E18').Value = '  +6.45%  '
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '331.82'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.75'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.46%  '
$ws.Range('E24').Value = '  +4.21%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('E28').Value = '  +6.74%  '
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '169.68'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.70'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.12%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  +5.70%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.420'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +11.66%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '39.45'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '314.91'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +9.76%  '
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '139.22'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0960'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0521'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.53'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.415'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +9.53%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.74'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('E51').Value = '  -0.15%  '
